$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ids = @("2021092003","2021092004","2021092005","2021092006","2021092007","2021092008","G2021092001","G2021092002")

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = 4 + $i

    # Column A: copy the number-format/border/font style from the row above
    # (style index 1) so the new numeric cell matches the existing rows,
    # then overwrite with the real value.
    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = 2 + $i

    # Column B: force a text entry (so "2021092003" etc. stay strings
    # instead of being parsed as numbers), then drop the number-format
    # override again so the cell ends up back on the default (unstyled)
    # look, matching the existing B2/B3 text cells.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $ids[$i]
    $ws.Cells.Item($row, 2).ClearFormats()
}
